$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 86.85714
$ws.Range("I4").Value = 86.85714
$ws.Range("K4").Value = 86.85714
$ws.Range("M4").Value = 27.14286
$ws.Range("H26").Value = 39962
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 39962
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 39962
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = -40650
$ws.Range("H112").Value = 1729.4333
$ws.Range("J112").Value = 1837.8077
$ws.Range("L112").Value = 5513.4231
$ws.Range("N112").Value = -7729.4231
$ws.Range("H121").Value = 1259.6595
$ws.Range("J121").Value = 1295.4222
$ws.Range("L121").Value = 3886.2666
$ws.Range("N121").Value = -7380.2666
$ws.Range("H129").Value = 969.9216
$ws.Range("I129").Value = 466.33334
$ws.Range("J129").Value = 1001.3958
$ws.Range("K129").Value = 1399.00002
$ws.Range("L129").Value = 3004.1874
$ws.Range("M129").Value = 3600.99998
$ws.Range("N129").Value = -13004.1874
$ws.Range("H131").Value = 7000
$ws.Range("J131").Value = 10000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("H137").Value = 3918.6938
$ws.Range("I137").Value = 3469.457
$ws.Range("J137").Value = 5041.7856
$ws.Range("K137").Value = 10408.371
$ws.Range("L137").Value = 15125.3568
$ws.Range("M137").Value = -7858.370999999999
$ws.Range("N137").Value = -20225.3568
$ws.Range("H138").Value = 4069.4487
$ws.Range("I138").Value = 2135.3333
$ws.Range("J138").Value = 4421.106
$ws.Range("K138").Value = 6405.999899999999
$ws.Range("L138").Value = 13263.318
$ws.Range("M138").Value = -1265.999899999999
$ws.Range("N138").Value = -23543.318
$ws.Range("H141").Value = 2427.5334
$ws.Range("I141").Value = 1913.04
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 5739.12
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -559.1199999999999
$ws.Range("N141").Value = -25360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8035.37
$ws.Range("I32").Value = 5991.027
$ws.Range("J32").Value = 14168.4
$ws.Range("K32").Value = 5991.027
$ws.Range("L32").Value = 14168.4
$ws.Range("M32").Value = -5704.027
$ws.Range("N32").Value = -14742.4
$ws.Range("H132").Value = 4206
$ws.Range("I132").Value = 2002
$ws.Range("J132").Value = 5969.2
$ws.Range("K132").Value = 6006
$ws.Range("L132").Value = 17907.6
$ws.Range("M132").Value = -3476
$ws.Range("N132").Value = -22967.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""
$ws.Range("H134").Value = 2909.0425
$ws.Range("I134").Value = 1427.8055
$ws.Range("J134").Value = 7756.727
$ws.Range("K134").Value = 4283.416499999999
$ws.Range("L134").Value = 23270.181
$ws.Range("M134").Value = -1748.416499999999
$ws.Range("N134").Value = -28340.181

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4325.5103
$ws.Range("I31").Value = 1782.7142
$ws.Range("J31").Value = 5342.6284
$ws.Range("K31").Value = 1782.7142
$ws.Range("L31").Value = 5342.6284
$ws.Range("M31").Value = -1487.7142
$ws.Range("N31").Value = -5932.6284
$ws.Range("H34").Value = 4325.5103
$ws.Range("I34").Value = 1782.7142
$ws.Range("J34").Value = 5342.6284
$ws.Range("K34").Value = 1782.7142
$ws.Range("L34").Value = 5342.6284
$ws.Range("M34").Value = -1580.7142
$ws.Range("N34").Value = -5746.6284
$ws.Range("H132").Value = 2316.6206
$ws.Range("I132").Value = 1775.74
$ws.Range("J132").Value = 5697.125
$ws.Range("K132").Value = 5327.22
$ws.Range("L132").Value = 17091.375
$ws.Range("M132").Value = -2797.22
$ws.Range("N132").Value = -22151.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2045
$ws.Range("I5").Value = 609.2857
$ws.Range("J5").Value = 4278.3335
$ws.Range("K5").Value = 1827.8571
$ws.Range("L5").Value = 12835.0005
$ws.Range("M5").Value = -1715.8571
$ws.Range("N5").Value = -13059.0005
$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -15354
$ws.Range("H113").Value = 617.1951
$ws.Range("I113").Value = 623.25
$ws.Range("K113").Value = 1869.75
$ws.Range("M113").Value = 300.25
$ws.Range("H135").Value = 2045
$ws.Range("I135").Value = 609.2857
$ws.Range("J135").Value = 4278.3335
$ws.Range("K135").Value = 5483.571300000001
$ws.Range("L135").Value = 38505.0015
$ws.Range("M135").Value = -2948.571300000001
$ws.Range("N135").Value = -43575.0015
$ws.Range("H137").Value = 7399.9287
$ws.Range("I137").Value = 2859.9443
$ws.Range("J137").Value = 15571.9
$ws.Range("K137").Value = 8579.832900000001
$ws.Range("L137").Value = 46715.7
$ws.Range("M137").Value = -3479.832900000001
$ws.Range("N137").Value = -56915.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1463.3334
$ws.Range("I113").Value = 1322
$ws.Range("J113").Value = 1640
$ws.Range("K113").Value = 1322
$ws.Range("L113").Value = 1640
$ws.Range("M113").Value = 848
$ws.Range("N113").Value = -5980
$ws.Range("H132").Value = 3722.4092
$ws.Range("I132").Value = 1867
$ws.Range("J132").Value = 4418.1875
$ws.Range("K132").Value = 5601
$ws.Range("L132").Value = 13254.5625
$ws.Range("M132").Value = -3071
$ws.Range("N132").Value = -18314.5625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4622.1113
$ws.Range("I7").Value = 2433.3333
$ws.Range("K7").Value = 2433.3333
$ws.Range("M7").Value = -2321.3333
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H100").Value = 1680.4375
$ws.Range("I100").Value = 1334.4546
$ws.Range("J100").Value = 2441.6
$ws.Range("K100").Value = 1334.4546
$ws.Range("L100").Value = 2441.6
$ws.Range("M100").Value = -793.4546
$ws.Range("N100").Value = -3523.6
$ws.Range("H126").Value = 4622.1113
$ws.Range("I126").Value = 2433.3333
$ws.Range("K126").Value = 7299.999899999999
$ws.Range("M126").Value = -4829.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17823.076
$ws.Range("J64").Value = 17823.076
$ws.Range("L64").Value = 17823.076
$ws.Range("N64").Value = -18319.076
$ws.Range("H67").Value = 17823.076
$ws.Range("J67").Value = 17823.076
$ws.Range("L67").Value = 17823.076
$ws.Range("N67").Value = -19539.076
$ws.Range("H113").Value = 9461.637000000001
$ws.Range("I113").Value = 14584.286
$ws.Range("J113").Value = 497
$ws.Range("K113").Value = 43752.858
$ws.Range("L113").Value = 1491
$ws.Range("M113").Value = -41582.858
$ws.Range("N113").Value = -5831
$ws.Range("H126").Value = 563283.5
$ws.Range("I126").Value = 2087.7778
$ws.Range("J126").Value = 1068359.8
$ws.Range("K126").Value = 6263.3334
$ws.Range("L126").Value = 3205079.4
$ws.Range("M126").Value = -3793.3334
$ws.Range("N126").Value = -3210019.4
$ws.Range("H132").Value = 4169348.2
$ws.Range("I132").Value = 2552.712
$ws.Range("J132").Value = 15876060
$ws.Range("K132").Value = 7658.136
$ws.Range("L132").Value = 47628180
$ws.Range("M132").Value = -5128.136
$ws.Range("N132").Value = -47633240
$ws.Range("H136").Value = 4777.2104
$ws.Range("I136").Value = 1723.4445
$ws.Range("J136").Value = 7525.6
$ws.Range("K136").Value = 5170.333500000001
$ws.Range("L136").Value = 22576.8
$ws.Range("M136").Value = -2620.333500000001
$ws.Range("N136").Value = -27676.8

Write-Host "Applied all Chocobo_Profits updates"